$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 73, shifting existing rows 73-100 down to 74-101.
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new weekly price entry.
$ws.Cells.Item(73, 1).Value = 11
$ws.Cells.Item(73, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(73, 3).Value = "Bíobío"
$ws.Cells.Item(73, 4).Value = 44468
$ws.Cells.Item(73, 5).Value = 8
$ws.Cells.Item(73, 6).Value = "Fruta"
$ws.Cells.Item(73, 7).Value = 100108
$ws.Cells.Item(73, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(73, 9).Value = 100108005
$ws.Cells.Item(73, 10).Value = "Piña"
$ws.Cells.Item(73, 11).Value = "Caramelo"
$ws.Cells.Item(73, 12).Value = "Segunda"
$ws.Cells.Item(73, 13).Value = 200
$ws.Cells.Item(73, 14).Value = 19000
$ws.Cells.Item(73, 15).Value = 20000
$ws.Cells.Item(73, 16).Value = 19500
$ws.Cells.Item(73, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(73, 18).Value = "Ecuador"
$ws.Cells.Item(73, 19).Value = 1393
$ws.Cells.Item(73, 20).Value = 14
